$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix row 8: swap B8/C8 values (ear_angle min/max was reversed)
$ws.Range("B8").Value = -50
$ws.Range("C8").Value = 50

# Header row for new columns D (name) and E (step)
$ws.Range("D1").Value = "name"
$ws.Range("E1").Value = "step"

# New column D ("pretty" display name) and E (step) values per row.
# Values are entered in this specific order so that the generated
# shared-strings table matches the source order.
$ws.Range("D2").Value = "Eye Aspect Ratio"
$ws.Range("D5").Value = "Face Aspect Ratio"
$ws.Range("D9").Value = "Ear Tip Angle"
$ws.Range("D3").Value = "Eye Distance"
$ws.Range("D4").Value = "Eye Height"
$ws.Range("D6").Value = "Nose Size"
$ws.Range("D7").Value = "Whisker Length"
$ws.Range("D8").Value = "Ear Angle"
$ws.Range("D10").Value = "Ear Point"
$ws.Range("D11").Value = "Ear Length"
$ws.Range("D12").Value = "Ear Orientation"
$ws.Range("D13").Value = "Fur Lightness"
$ws.Range("D14").Value = "Fur Saturation"

$ws.Range("E2").Value = 0.01
$ws.Range("E3").Value = 1
$ws.Range("E4").Value = 1
$ws.Range("E5").Value = 0.01
$ws.Range("E6").Value = 1
$ws.Range("E7").Value = 1
$ws.Range("E8").Value = 0.1
$ws.Range("E9").Value = 0.1
$ws.Range("E10").Value = 0.1
$ws.Range("E11").Value = 1
$ws.Range("E12").Value = 0.01
$ws.Range("E13").Value = 0.1
$ws.Range("E14").Value = 0.1

# Column D width (COM ColumnWidth has a constant +5/6 offset vs. the
# stored OOXML character width, so back it out to land exactly on 17)
$ws.Columns.Item(4).ColumnWidth = 16.166666666666668

# Update selection to match the saved view state
$ws.Range("B10").Select()
